$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E (rows 2-51) to Text format so that numeric-looking
# strings (e.g. "0.739", "40.196.25") are preserved as text, matching the
# inlineStr cell type used in the source data.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '40.196.25'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '2.227.13'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '293.78'
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("D6").Value = '88.30'
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("D7").Value = '0.515'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '30.81'
$ws.Range("E10").Value = '  +0.43%  '
$ws.Range("D11").Value = '50.85'
$ws.Range("E11").Value = '  +6.41%  '
$ws.Range("D12").Value = '0.0784'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("E13").Value = '  +3.35%  '
$ws.Range("D14").Value = '6.45'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.544.90'
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '13.86'
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.206.33'
$ws.Range("E17").Value = '  -0.51%  '
$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").Value = '0.739'
$ws.Range("E18").Value = '  +1.28%  '
$ws.Range("D19").Value = '40.129.15'
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D21").Value = '11.31'
$ws.Range("E21").Value = '  -3.49%  '
$ws.Range("D22").Value = '5.80'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").Value = '65.74'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = '236.32'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  +1.14%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '23.33'
$ws.Range("E28").Value = '  +3.24%  '
$ws.Range("E29").Value = '  +1.33%  '
$ws.Range("D30").Value = '2.07'
$ws.Range("E30").Value = '  -10.06%  '
$ws.Range("D31").Value = '159.16'
$ws.Range("E31").Value = '  +3.97%  '
$ws.Range("D32").Value = '32.02'
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").Value = '3.02'
$ws.Range("E35").Value = '  +6.64%  '
$ws.Range("D36").Value = '0.0718'
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("E38").Value = '  +1.78%  '
$ws.Range("D39").Value = '1.77'
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("E40").Value = '  -0.29%  '
$ws.Range("D41").Value = '15.75'
$ws.Range("E41").Value = '  -0.60%  '
$ws.Range("D42").Value = '2.090.62'
$ws.Range("E42").Value = '  -0.77%  '
$ws.Range("D43").Value = '3.77'
$ws.Range("E43").Value = '  -2.46%  '
$ws.Range("D44").Value = '19.43'
$ws.Range("E44").Value = '  +10.44%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0271'
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '10.10'
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").Value = '2.77'
$ws.Range("E47").Value = '  +3.40%  '
$ws.Range("E48").Value = '  -13.39%  '
$ws.Range("D49").Value = '2.432.59'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '1.48'
$ws.Range("E50").Value = '  +1.63%  '
$ws.Range("E51").Value = '  +3.73%  '

# Restore the default "Normal" style so no stray number-format styling is left behind.
$dataRange.Style = "Normal"
